$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 2 4 "63.815.73"
Set-TextValue 2 5 "  -0.11%  "
Set-TextValue 3 4 "2.749.73"
Set-TextValue 3 5 "  -0.11%  "
Set-TextValue 4 5 "  -0.18%  "
Set-TextValue 5 4 "573.28"
Set-TextValue 5 5 "  -1.29%  "
Set-TextValue 6 4 "157.46"
Set-TextValue 6 5 "  +0.81%  "
Set-TextValue 7 5 "  +0.06%  "
Set-TextValue 8 5 "  -1.53%  "
Set-TextValue 9 5 "  -3.25%  "
Set-TextValue 10 5 "  -0.02%  "
Set-TextValue 11 2 "Toncoin"
Set-TextValue 11 3 "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue 11 4 "5.64"
Set-TextValue 11 5 "  -16.63%  "
Set-TextValue 12 2 "Cardano"
Set-TextValue 12 3 "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue 12 4 "0.381"
Set-TextValue 12 5 "  -2.48%  "
Set-TextValue 13 4 "3.234.84"
Set-TextValue 13 5 "  -0.77%  "
Set-TextValue 14 4 "26.46"
Set-TextValue 14 5 "  -1.58%  "
Set-TextValue 15 4 "63.518.92"
Set-TextValue 15 5 "  -0.56%  "
Set-TextValue 16 5 "  -2.65%  "
Set-TextValue 17 4 "2.750.99"
Set-TextValue 17 5 "  -0.74%  "
Set-TextValue 18 4 "12.13"
Set-TextValue 18 5 "  +0.89%  "
Set-TextValue 19 5 "  -2.15%  "
Set-TextValue 20 4 "354.67"
Set-TextValue 21 5 "  -4.01%  "
Set-TextValue 22 5 "  +0.89%  "
Set-TextValue 23 4 "0.998"
Set-TextValue 23 5 "  -0.23%  "
Set-TextValue 24 4 "65.12"
Set-TextValue 24 5 "  -2.08%  "
Set-TextValue 25 4 "0.170"
Set-TextValue 25 5 "  -0.92%  "
Set-TextValue 26 5 "  +0.00%  "
Set-TextValue 27 5 "  -2.30%  "
Set-TextValue 28 4 "0.0₃0902"
Set-TextValue 28 5 "  -0.44%  "
Set-TextValue 29 5 "  -4.10%  "
Set-TextValue 30 4 "6.96"
Set-TextValue 30 5 "  -2.53%  "
Set-TextValue 31 4 "169.09"
Set-TextValue 31 5 "  -2.83%  "
Set-TextValue 32 4 "1.19"
Set-TextValue 32 5 "  -7.22%  "
Set-TextValue 33 4 "20.12"
Set-TextValue 33 5 "  -2.12%  "
Set-TextValue 35 4 "4.85"
Set-TextValue 35 5 "  -0.65%  "
Set-TextValue 36 5 "  -0.98%  "
Set-TextValue 37 5 "  -2.68%  "
Set-TextValue 38 5 "  -3.77%  "
Set-TextValue 39 4 "6.17"
Set-TextValue 39 5 "  +5.66%  "
Set-TextValue 40 5 "  -2.99%  "
Set-TextValue 41 4 "325.76"
Set-TextValue 41 5 "  -5.24%  "
Set-TextValue 42 4 "38.81"
Set-TextValue 42 5 "  -1.28%  "
Set-TextValue 43 4 "21.35"
Set-TextValue 43 5 "  -2.99%  "
Set-TextValue 44 5 "  -1.97%  "
Set-TextValue 45 4 "21.28"
Set-TextValue 45 5 "  -3.42%  "
Set-TextValue 46 5 "  -1.41%  "
Set-TextValue 47 4 "135.10"
Set-TextValue 47 5 "  -2.19%  "
Set-TextValue 48 2 "Mantle"
Set-TextValue 48 3 "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue 48 4 "0.624"
Set-TextValue 48 5 "  -3.98%  "
Set-TextValue 49 2 "Stellar"
Set-TextValue 49 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue 49 4 "0.101"
Set-TextValue 49 5 "  -0.95%  "
Set-TextValue 50 5 "  +0.37%  "
Set-TextValue 51 5 "  +0.33%  "
